$wb = $excel.ActiveWorkbook

# ---- Sheet: Summary ----
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B2").Value = 0.2740213523131673
$ws.Range("C2").Value = 0.06422018348623854
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0.1206896551724138
$ws.Range("F2").Value = 0.2554744525547445
$ws.Range("G2").Value = 0.6408450704225352
$ws.Range("H2").Value = 0.8127340823970037
$ws.Range("I2").Value = 28
$ws.Range("J2").Value = 408
$ws.Range("K2").Value = 126
$ws.Range("L2").Value = 0

# ---- Sheet: Classification Report ----
$ws = $wb.Worksheets.Item("Classification Report")

# row 2 -> class "0"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0.2359550561797753
$ws.Range("D2").Value = 0.3818181818181818

# row 3 -> class "1"
$ws.Range("B3").Value = 0.06422018348623854
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0.1206896551724138

# row 4 -> accuracy
$ws.Range("B4").Value = 0.2740213523131673
$ws.Range("C4").Value = 0.2740213523131673
$ws.Range("D4").Value = 0.2740213523131673
$ws.Range("E4").Value = 0.2740213523131673

# row 5 -> macro avg
$ws.Range("B5").Value = 0.5321100917431193
$ws.Range("C5").Value = 0.6179775280898876
$ws.Range("D5").Value = 0.2512539184952978

# row 6 -> weighted avg
$ws.Range("B6").Value = 0.9533775180384603
$ws.Range("C6").Value = 0.2740213523131673
$ws.Range("D6").Value = 0.3688082196365421

# ---- Sheet: Confusion Matrix ----
$ws = $wb.Worksheets.Item("Confusion Matrix")

# row 2 -> Actual 0
$ws.Range("B2").Value = 126
$ws.Range("C2").Value = 408

# row 3 -> Actual 1
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 28
